# Courts.xlsx - Add files via upload
# On the "CourtReports" sheet, insert a new column before column E
# (NEW_BTN / "Click") shifting the existing REPORT_TYPE..CREATED_BY_VERIFY
# columns one place to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("CourtReports")
# ("CourtReports" is also the workbook's ActiveSheet / selected tab.)

# Insert a new column at E; Excel copies formatting from the column to
# the left (D) by default, matching the observed column-width/style merge.
$dWidth = $ws.Columns("D:D").ColumnWidth
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = $dWidth

# Populate the newly inserted column E.
$ws.Range("E1").Value = "NEW_BTN"
$ws.Range("E2").Value = "Click"

# Update the active selection to A3, as recorded in the saved view state.
$ws.Range("A3").Select()
